$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Step 1: row 1 ("100" -> "0M") ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"

# --- Step 2: insert 12 new rows right after row 1 (before original row 2) ---
$block1 = @("0M", "0M", "20", "0.00002", "0.00006", "0.00003", "0.00001", "0.00003", "0.00003", "0.00004", "0.00069", "100.0")
$insertIndex = 2
foreach ($val in $block1) {
    $beforeRow = $t.Rows.Item($insertIndex)
    $newRow = $t.Rows.Add($beforeRow)
    $newRow.Cells.Item(1).Range.Text = $val
    $insertIndex = $insertIndex + 1
}

# After step 2: original row 2 ("0") has shifted to index 14 (1 + 12 + 1)
# original row 3 ("70") is now at index 15

# --- Step 3: row 15 ("70" -> "0.00000") ---
$t.Rows.Item(15).Cells.Item(1).Range.Text = "0.00000"

# --- Step 4: insert 9 new rows right after row 15 (before original row 4, now at index 16) ---
$block2 = @("0.00000", "0.00000", "0.00000", "0.00000", "0.00000", "0.00000", "0.00000", "0.00000", "0.0")
$insertIndex = 16
foreach ($val in $block2) {
    $beforeRow = $t.Rows.Item($insertIndex)
    $newRow = $t.Rows.Add($beforeRow)
    $newRow.Cells.Item(1).Range.Text = $val
    $insertIndex = $insertIndex + 1
}

# After step 4: table has 47 rows.
# original row 4 ("0") now at index 25
# original rows 4-23 (20 rows, unchanged) now occupy indices 25-44
# original row 24 (multi-value "20...100.0") now at index 45 -> DELETE
# original row 25 (empty) now at index 45 (after delete) -> set text "100"
# original row 26 (multi-value "0...0") now at index 46 -> set text "0"

# --- Step 5: delete original row 24 (now at index 45) ---
$t.Rows.Item(45).Delete()

# --- Step 6: original row 25 (now at index 45) gets text "100" ---
$t.Rows.Item(45).Cells.Item(1).Range.Text = "100"

# --- Step 7: original row 26 (now at index 46) gets text "0" ---
$t.Rows.Item(46).Cells.Item(1).Range.Text = "0"

# --- Step 8: append a new row with "70" at the end ---
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "70"

Write-Host "Final row count:" $t.Rows.Count
